$d = $word.ActiveDocument

function Set-ParaText($para, $newText) {
    $r = $d.Range($para.Range.Start, $para.Range.End - 1)
    $r.Text = $newText
}

function Expect-Contains($para, $needle) {
    if ($para.Range.Text.IndexOf($needle) -lt 0) {
        throw "Paragraph text [$($para.Range.Text)] did not contain expected [$needle]"
    }
}

# Locate the "Suggestions for improvements" heading paragraph so we only touch
# that section below it (the same sentences also appear earlier in the doc,
# verbatim, in the "Not clear:" / "Review Completeness" sections, and those
# must stay untouched).
$sugIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.TrimEnd() -eq "Suggestions for improvements") {
        $sugIndex = $i
        break
    }
}
if ($sugIndex -eq -1) {
    throw "Could not find 'Suggestions for improvements' heading"
}

# Walk the bullet paragraphs that follow the heading. They alternate with
# blank paragraphs: heading, bullet, blank, bullet, blank, ...
$idx = $sugIndex + 1

# 1) "- Spelling/grammar error: ..." -> add "Fix the " before the phrase.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "Spelling/grammar error: "
Set-ParaText $para "- Fix the Spelling/grammar error: Association between BoatType and Boat. Actual: Is of an, Expected: Is of a."
$idx += 2

# 2) "- Switch the right column..." paragraph is unchanged, skip over it.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "Switch the right column"
$idx += 2

# 3) "- I could not figure out..." paragraph gets rewritten, and the content
#    that used to start the *next* paragraph ("Alternatively keep...") is
#    merged into it; the rest of that next paragraph becomes its own new
#    bullet ("- Rename the association...").
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "I could not figure out"
Set-ParaText $para "- Remove the Reservation class or alternatively keep the Reservation, but remove the association between Secretary and Berth, since the Secretary approving the Reservation automatically assigns Berth."
$idx += 2

$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "Alternatively keep the Reservation"
Set-ParaText $para "- Rename the association between Reservation and Berth from reserves to assigns."
$idx += 2

# 4) "- Both Member and Secretary includes..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "Both Member and Secretary"
Set-ParaText $para "- Make a separate class (maybe Person?) for the attributes username and password. This information is repeated in both Member and Secretary classes and could therefore be put in a separate class for the two to make use of (Larman, chapter 9, figure 9.9, 2004)."
$idx += 2

# 5) "- Missing a class for membership fee..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "Missing a class for membership fee"
Set-ParaText $para "- Create a class for MembershipFee with attributes such as fixedCost and variableCost."
$idx += 2

# 6) "- The classes Boat and Berth should have an association..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "The classes Boat and Berth should have an association"
Set-ParaText $para "- Create an association between the classes Boat and Berth."
$idx += 2

# 7) "- In the class Boat the attribute Picture should be made optional..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "In the class Boat the attribute Pic"
Set-ParaText $para "- Make the attribute Picture in the class Boat optional [0..1] (Larman, chapter 9, figure 9.20, 2004)."
$idx += 2

# 8) "- The class Berth includes contradicting attributes..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "The class Berth includes contradicting attributes"
Set-ParaText $para "- Choose one of the attributes IsAvailable or IsReserved in the class Berth, since they are now contradicting."
$idx += 2

# 9) "- The class Calendar might not need the attribute Title..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "The class Calendar might not need the attribute Title"
Set-ParaText $para "- Remove the attribute Title in the class Calendar, since there is only one Boat Club Calendar."
$idx += 2

# 10) "- The requirements state that the Secretary manages events..." paragraph.
$para = $d.Paragraphs.Item($idx)
Expect-Contains $para "The requirements state that the Secretary manages events"
Set-ParaText $para "- Switch places between Calendar and Event in the domain model, since the requirements state that the Secretary manages events."
$idx += 2

Write-Host "Suggestions section rewritten (heading at paragraph $sugIndex)."
